# Update the "Date" column (column B) values for the test-run rows on
# Sheet1, adding a fresh batch of execution timestamps for the newly added
# "Personal_SearchTransaction_Generic_TC" test data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(
    "Thu Mar 07 23:54:01 EST 2024",
    "Thu Mar 07 23:54:46 EST 2024",
    "Thu Mar 07 23:55:06 EST 2024",
    "Thu Mar 07 23:55:48 EST 2024",
    "Thu Mar 07 23:56:13 EST 2024",
    "Thu Mar 07 23:56:58 EST 2024",
    "Thu Mar 07 23:57:19 EST 2024",
    "Thu Mar 07 23:57:44 EST 2024",
    "Thu Mar 07 23:58:14 EST 2024",
    "Thu Mar 07 23:58:44 EST 2024",
    "Thu Mar 07 23:59:19 EST 2024",
    "Thu Mar 07 23:59:41 EST 2024",
    "Thu Mar 07 23:59:56 EST 2024",
    "Fri Mar 08 00:00:24 EST 2024",
    "Fri Mar 08 00:01:07 EST 2024",
    "Fri Mar 08 00:01:25 EST 2024",
    "Fri Mar 08 00:01:57 EST 2024",
    "Fri Mar 08 00:02:36 EST 2024",
    "Fri Mar 08 00:03:06 EST 2024",
    "Fri Mar 08 00:03:31 EST 2024",
    "Fri Mar 08 00:04:12 EST 2024",
    "Fri Mar 08 00:05:00 EST 2024",
    "Fri Mar 08 00:05:20 EST 2024",
    "Fri Mar 08 00:05:54 EST 2024",
    "Fri Mar 08 00:06:30 EST 2024",
    "Fri Mar 08 00:06:57 EST 2024",
    "Fri Mar 08 00:07:31 EST 2024",
    "Fri Mar 08 00:07:46 EST 2024",
    "Fri Mar 08 00:08:16 EST 2024",
    "Fri Mar 08 00:08:42 EST 2024",
    "Fri Mar 08 00:09:15 EST 2024",
    "Fri Mar 08 00:09:50 EST 2024",
    "Fri Mar 08 00:10:23 EST 2024",
    "Fri Mar 08 00:10:47 EST 2024",
    "Fri Mar 08 00:11:26 EST 2024",
    "Fri Mar 08 00:12:00 EST 2024",
    "Fri Mar 08 00:12:31 EST 2024",
    "Fri Mar 08 00:13:02 EST 2024",
    "Fri Mar 08 00:13:45 EST 2024",
    "Fri Mar 08 00:14:26 EST 2024",
    "Fri Mar 08 00:14:57 EST 2024"
)

# Rows 2-18 and 25-48 get the new run timestamps, in order; rows 19-24 are
# untouched (they hold a different, unrelated set of values in column B).
$rows = @(2..18) + @(25..48)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $dates[$i]
}
